$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-9
# from serial date 45185 (2023-09-16) to 45204 (2023-10-05),
# preserving existing cell formatting.
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
